$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.022.67"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.640.19"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.55%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.62%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.79"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5157"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.49%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.61%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2587"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06380"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.66%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.87"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.62%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07767"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.300"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.647.77"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.20%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5481"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0₅7782"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.62"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.98%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.040.83"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("E18").Value = "  -0.68%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "199.24"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.466"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.990"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.118"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.891"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.28%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "142.29"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1233"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +7.56%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "6.880"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.63"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.64%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.243"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.04865"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.75%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.311"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.94%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.243"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.542"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.379"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9205"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.5598"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.570"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.113.53"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01574"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.533"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.10%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.573"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.83%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8092"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.01%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "99.56"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("E45").Value = "  -0.05%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.780.29"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4537"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "55.35"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05223"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.09596"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
